# Updated CVDs for the month
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Grafton Wisconsin  (Professional Voluntary Turnover ytd + Commit/Forecast
# monthly CVD values; Internal Fill Rate O7 cleared)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Grafton Wisconsin")
$ws.Range("E2").Value = 0.0815
$ws.Range("E3").Value = 0.0815
$ws.Range("E4").Value = 0.0815
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("O7").ClearContents()

# ---------------------------------------------------------------------------
# Guadalajara Mexico  (Commit/Forecast monthly CVD values)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Guadalajara Mexico")
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# ---------------------------------------------------------------------------
# Hyderabad India  (Professional Voluntary Turnover ytd + Commit/Forecast
# monthly CVD values)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyderabad India")
$ws.Range("E2").Value = 0.0303
$ws.Range("E3").Value = 0.0303
$ws.Range("E4").Value = 0.0303
$ws.Range("I4").Value = 0.0154
$ws.Range("J4").Value = 0.0153
$ws.Range("L4").Value = 0.0154
$ws.Range("N4").Value = 0.0152
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

# ---------------------------------------------------------------------------
# Black River Falls Wisconsin  (Professional Voluntary Turnover ytd +
# Commit/Forecast monthly CVD values)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Black River Falls Wisconsin")
$ws.Range("E2").Value = 0.7143
$ws.Range("E3").Value = 0.7143
$ws.Range("E4").Value = 0.7143
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# ---------------------------------------------------------------------------
# South Beloit Gardner St Illino  (Commit/Forecast monthly CVD values)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("South Beloit Gardner St Illino")
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

# ---------------------------------------------------------------------------
# Chicago Lasalle Illinois  (Commit/Forecast monthly CVD values)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Chicago Lasalle Illinois")
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1
